$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.039483
$ws.Range("H2").Value = 3.118449
$ws.Range("I2").Value = 0.01443400247257094
$ws.Range("J2").Value = 0.01491631751620889
$ws.Range("M2").Value = 19.77408333333333
$ws.Range("N2").Value = 59.32225
$ws.Range("O2").Value = 0.3380388258879848
$ws.Range("P2").Value = 0.339186328349942
$ws.Range("Q2").Value = 20.55482346558333
$ws.Range("R2").Value = 184.99341119025
$ws.Range("S2").Value = 0.004879253248692151
$ws.Range("T2").Value = 0.00505941097082482
$ws.Range("G3").Value = 1.039483
$ws.Range("H3").Value = 3.118449
$ws.Range("I3").Value = 0.01443400247257094
$ws.Range("J3").Value = 0.01491631751620889
$ws.Range("O3").Value = 0.3069959581674471
$ws.Range("P3").Value = 0.3080380828904952
$ws.Range("Q3").Value = 18.66722767186067
$ws.Range("R3").Value = 168.005049046746
$ws.Range("S3").Value = 0.004431180419258217
$ws.Range("T3").Value = 0.0045947938514789
$ws.Range("G4").Value = 1.039483
$ws.Range("H4").Value = 3.118449
$ws.Range("I4").Value = 0.01443400247257094
$ws.Range("J4").Value = 0.01491631751620889
$ws.Range("M4").Value = 12.46730333333333
$ws.Range("N4").Value = 37.40191
$ws.Range("O4").Value = 0.2131290998296268
$ws.Range("P4").Value = 0.2138525852639604
$ws.Range("Q4").Value = 12.95954987084333
$ws.Range("R4").Value = 116.63594883759
$ws.Range("S4").Value = 0.003076305953917653
$ws.Range("T4").Value = 0.003189893063459368
$ws.Range("G5").Value = 1.039483
$ws.Range("H5").Value = 3.118449
$ws.Range("I5").Value = 0.01443400247257094
$ws.Range("J5").Value = 0.01491631751620889
$ws.Range("M5").Value = 0.593699
$ws.Range("N5").Value = 1.187398
$ws.Range("O5").Value = 0.01014931056513554
$ws.Range("P5").Value = 0.006789175527058808
$ws.Range("Q5").Value = 0.617140017617
$ws.Range("R5").Value = 3.702840105702
$ws.Range("S5").Value = 0.0001464951737920568
$ws.Range("T5").Value = 0.000101269497834884
$ws.Range("G6").Value = 1.039483
$ws.Range("H6").Value = 3.118449
$ws.Range("I6").Value = 0.01443400247257094
$ws.Range("J6").Value = 0.01491631751620889
$ws.Range("M6").Value = 7.703215333333333
$ws.Range("N6").Value = 23.109646
$ws.Range("O6").Value = 0.1316868055498057
$ws.Range("P6").Value = 0.1321338279685434
$ws.Range("Q6").Value = 8.007361384339331
$ws.Range("R6").Value = 72.066252459054
$ws.Range("S6").Value = 0.001900767676910865
$ws.Range("T6").Value = 0.001970950132610916
$ws.Range("I7").Value = 0.8791289547788569
$ws.Range("J7").Value = 0.9085052224491242
$ws.Range("M7").Value = 19.77408333333333
$ws.Range("N7").Value = 59.32225
$ws.Range("O7").Value = 0.3380388258879848
$ws.Range("P7").Value = 0.339186328349942
$ws.Range("Q7").Value = 1251.928597303584
$ws.Range("R7").Value = 11267.35737573225
$ws.Range("S7").Value = 0.2971797196775761
$ws.Range("T7").Value = 0.3081525506892658
$ws.Range("I8").Value = 0.8791289547788569
$ws.Range("J8").Value = 0.9085052224491242
$ws.Range("O8").Value = 0.3069959581674471
$ws.Range("P8").Value = 0.3080380828904952
$ws.Range("S8").Value = 0.2698890358250814
$ws.Range("T8").Value = 0.2798542070192311
$ws.Range("I9").Value = 0.8791289547788569
$ws.Range("J9").Value = 0.9085052224491242
$ws.Range("M9").Value = 12.46730333333333
$ws.Range("N9").Value = 37.40191
$ws.Range("O9").Value = 0.2131290998296268
$ws.Range("P9").Value = 0.2138525852639604
$ws.Range("Q9").Value = 789.3247596437234
$ws.Range("R9").Value = 7103.92283679351
$ws.Range("S9").Value = 0.1873679627661785
$ws.Range("T9").Value = 0.1942861905465547
$ws.Range("I10").Value = 0.8791289547788569
$ws.Range("J10").Value = 0.9085052224491242
$ws.Range("M10").Value = 0.593699
$ws.Range("N10").Value = 1.187398
$ws.Range("O10").Value = 0.01014931056513554
$ws.Range("P10").Value = 0.006789175527058808
$ws.Range("Q10").Value = 37.588025890313
$ws.Range("R10").Value = 225.528155341878
$ws.Range("S10").Value = 0.00892255278885362
$ws.Range("T10").Value = 0.006168001422456712
$ws.Range("I11").Value = 0.8791289547788569
$ws.Range("J11").Value = 0.9085052224491242
$ws.Range("M11").Value = 7.703215333333333
$ws.Range("N11").Value = 23.109646
$ws.Range("O11").Value = 0.1316868055498057
$ws.Range("P11").Value = 0.1321338279685434
$ws.Range("Q11").Value = 487.7027877560673
$ws.Range("R11").Value = 4389.325089804606
$ws.Range("S11").Value = 0.1157696837211673
$ws.Range("T11").Value = 0.1200442727716158
$ws.Range("G12").Value = 0.3690693333333333
$ws.Range("H12").Value = 1.107208
$ws.Range("I12").Value = 0.005124804994293743
$ws.Range("J12").Value = 0.005296051365434103
$ws.Range("M12").Value = 19.77408333333333
$ws.Range("N12").Value = 59.32225
$ws.Range("O12").Value = 0.3380388258879848
$ws.Range("P12").Value = 0.339186328349942
$ws.Range("Q12").Value = 7.29800775311111
$ws.Range("R12").Value = 65.682069778
$ws.Range("S12").Value = 0.001732383063175937
$ws.Range("T12").Value = 0.00179634821739429
$ws.Range("G13").Value = 0.3690693333333333
$ws.Range("H13").Value = 1.107208
$ws.Range("I13").Value = 0.005124804994293743
$ws.Range("J13").Value = 0.005296051365434103
$ws.Range("O13").Value = 0.3069959581674471
$ws.Range("P13").Value = 0.3080380828904952
$ws.Range("Q13").Value = 6.627815242803555
$ws.Range("R13").Value = 59.650337185232
$ws.Range("S13").Value = 0.001573294419644526
$ws.Range("T13").Value = 0.001631385509497911
$ws.Range("G14").Value = 0.3690693333333333
$ws.Range("H14").Value = 1.107208
$ws.Range("I14").Value = 0.005124804994293743
$ws.Range("J14").Value = 0.005296051365434103
$ws.Range("M14").Value = 12.46730333333333
$ws.Range("N14").Value = 37.40191
$ws.Range("O14").Value = 0.2131290998296268
$ws.Range("P14").Value = 0.2138525852639604
$ws.Range("Q14").Value = 4.601299329697778
$ws.Range("R14").Value = 41.41169396728
$ws.Range("S14").Value = 0.001092245075236201
$ws.Range("T14").Value = 0.001132574276188811
$ws.Range("G15").Value = 0.3690693333333333
$ws.Range("H15").Value = 1.107208
$ws.Range("I15").Value = 0.005124804994293743
$ws.Range("J15").Value = 0.005296051365434103
$ws.Range("M15").Value = 0.593699
$ws.Range("N15").Value = 1.187398
$ws.Range("O15").Value = 0.01014931056513554
$ws.Range("P15").Value = 0.006789175527058808
$ws.Range("Q15").Value = 0.2191160941306666
$ws.Range("R15").Value = 1.314696564784
$ws.Range("S15").Value = 0.00005201323747284489
$ws.Range("T15").Value = 0.00003595582232025159
$ws.Range("G16").Value = 0.3690693333333333
$ws.Range("H16").Value = 1.107208
$ws.Range("I16").Value = 0.005124804994293743
$ws.Range("J16").Value = 0.005296051365434103
$ws.Range("M16").Value = 7.703215333333333
$ws.Range("N16").Value = 23.109646
$ws.Range("O16").Value = 0.1316868055498057
$ws.Range("P16").Value = 0.1321338279685434
$ws.Range("Q16").Value = 2.843020547596444
$ws.Range("R16").Value = 25.587184928368
$ws.Range("S16").Value = 0.0006748691987642335
$ws.Range("T16").Value = 0.0006997875400328389
$ws.Range("G17").Value = 6.985879
$ws.Range("H17").Value = 13.971758
$ws.Range("I17").Value = 0.0970041787687547
$ws.Range("J17").Value = 0.06683039504177611
$ws.Range("M17").Value = 19.77408333333333
$ws.Range("N17").Value = 59.32225
$ws.Range("O17").Value = 0.3380388258879848
$ws.Range("P17").Value = 0.339186328349942
$ws.Range("Q17").Value = 138.1393535025833
$ws.Range("R17").Value = 828.8361210155
$ws.Range("S17").Value = 0.03279117869721802
$ws.Range("T17").Value = 0.02266795631639621
$ws.Range("G18").Value = 6.985879
$ws.Range("H18").Value = 13.971758
$ws.Range("I18").Value = 0.0970041787687547
$ws.Range("J18").Value = 0.06683039504177611
$ws.Range("O18").Value = 0.3069959581674471
$ws.Range("P18").Value = 0.3080380828904952
$ws.Range("Q18").Value = 125.4537051409887
$ws.Range("R18").Value = 752.722230845932
$ws.Range("S18").Value = 0.02977989080736018
$ws.Range("T18").Value = 0.02058630676748317
$ws.Range("G19").Value = 6.985879
$ws.Range("H19").Value = 13.971758
$ws.Range("I19").Value = 0.0970041787687547
$ws.Range("J19").Value = 0.06683039504177611
$ws.Range("M19").Value = 12.46730333333333
$ws.Range("N19").Value = 37.40191
$ws.Range("O19").Value = 0.2131290998296268
$ws.Range("P19").Value = 0.2138525852639604
$ws.Range("Q19").Value = 87.09507254296334
$ws.Range("R19").Value = 522.57043525778
$ws.Range("S19").Value = 0.02067441330069688
$ws.Range("T19").Value = 0.01429185275389559
$ws.Range("G20").Value = 6.985879
$ws.Range("H20").Value = 13.971758
$ws.Range("I20").Value = 0.0970041787687547
$ws.Range("J20").Value = 0.06683039504177611
$ws.Range("M20").Value = 0.593699
$ws.Range("N20").Value = 1.187398
$ws.Range("O20").Value = 0.01014931056513554
$ws.Range("P20").Value = 0.006789175527058808
$ws.Range("Q20").Value = 4.147509376421
$ws.Range("R20").Value = 16.590037505684
$ws.Range("S20").Value = 0.0009845255364400191
$ws.Range("T20").Value = 0.0004537232824812987
$ws.Range("G21").Value = 6.985879
$ws.Range("H21").Value = 13.971758
$ws.Range("I21").Value = 0.0970041787687547
$ws.Range("J21").Value = 0.06683039504177611
$ws.Range("M21").Value = 7.703215333333333
$ws.Range("N21").Value = 23.109646
$ws.Range("O21").Value = 0.1316868055498057
$ws.Range("P21").Value = 0.1321338279685434
$ws.Range("Q21").Value = 53.81373022961132
$ws.Range("R21").Value = 322.882381377668
$ws.Range("S21").Value = 0.01277417042703959
$ws.Range("T21").Value = 0.008830555921519838
$ws.Range("G22").Value = 0.3102503333333334
$ws.Range("H22").Value = 0.930751
$ws.Range("I22").Value = 0.004308058985523854
$ws.Range("J22").Value = 0.004452013627456771
$ws.Range("M22").Value = 19.77408333333333
$ws.Range("N22").Value = 59.32225
$ws.Range("O22").Value = 0.3380388258879848
$ws.Range("P22").Value = 0.339186328349942
$ws.Range("Q22").Value = 6.134915945527778
$ws.Range("R22").Value = 55.21424350975001
$ws.Range("S22").Value = 0.001456291201322667
$ws.Range("T22").Value = 0.001510062156060969
$ws.Range("G23").Value = 0.3102503333333334
$ws.Range("H23").Value = 0.930751
$ws.Range("I23").Value = 0.004308058985523854
$ws.Range("J23").Value = 0.004452013627456771
$ws.Range("O23").Value = 0.3069959581674471
$ws.Range("P23").Value = 0.3080380828904952
$ws.Range("Q23").Value = 5.571532778894889
$ws.Range("R23").Value = 50.14379501005401
$ws.Range("S23").Value = 0.001322556696102776
$ws.Range("T23").Value = 0.001371389742804143
$ws.Range("G24").Value = 0.3102503333333334
$ws.Range("H24").Value = 0.930751
$ws.Range("I24").Value = 0.004308058985523854
$ws.Range("J24").Value = 0.004452013627456771
$ws.Range("M24").Value = 12.46730333333333
$ws.Range("N24").Value = 37.40191
$ws.Range("O24").Value = 0.2131290998296268
$ws.Range("P24").Value = 0.2138525852639604
$ws.Range("Q24").Value = 3.867985014934445
$ws.Range("R24").Value = 34.81186513441
$ws.Range("S24").Value = 0.0009181727335976343
$ws.Range("T24").Value = 0.000952074623862013
$ws.Range("G25").Value = 0.3102503333333334
$ws.Range("H25").Value = 0.930751
$ws.Range("I25").Value = 0.004308058985523854
$ws.Range("J25").Value = 0.004452013627456771
$ws.Range("M25").Value = 0.593699
$ws.Range("N25").Value = 1.187398
$ws.Range("O25").Value = 0.01014931056513554
$ws.Range("P25").Value = 0.006789175527058808
$ws.Range("Q25").Value = 0.1841953126496667
$ws.Range("R25").Value = 1.105171875898
$ws.Range("S25").Value = 0.00004372382857700437
$ws.Range("T25").Value = 0.00003022550196566182
$ws.Range("G26").Value = 0.3102503333333334
$ws.Range("H26").Value = 0.930751
$ws.Range("I26").Value = 0.004308058985523854
$ws.Range("J26").Value = 0.004452013627456771
$ws.Range("M26").Value = 7.703215333333333
$ws.Range("N26").Value = 23.109646
$ws.Range("O26").Value = 0.1316868055498057
$ws.Range("P26").Value = 0.1321338279685434
$ws.Range("Q26").Value = 2.389925124905111
$ws.Range("R26").Value = 21.509326124146
$ws.Range("S26").Value = 0.0005673145259237733
$ws.Range("T26").Value = 0.0005882616027639837
